$d = $word.ActiveDocument

# Locate the four target paragraphs (TEAM_LEADER, INSPECTOR, INSPECTOR_2,
# INSPECTOR_3) inside the "TO" table cell by scanning for their known
# placeholder text. This is more robust than hard-coded paragraph indices.
$paraTeamLeader = $null
$paraInspector1 = $null
$paraInspector2 = $null
$paraInspector3 = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*{TEAM_LEADER}*") {
        $paraTeamLeader = $p
    } elseif ($t -like "*{INSPECTOR_2}*") {
        $paraInspector2 = $p
    } elseif ($t -like "*{INSPECTOR_3}*") {
        $paraInspector3 = $p
    } elseif ($t -like "*{INSPECTOR}*") {
        $paraInspector1 = $p
    }
}

# Paragraph.Range alone does not reliably clip Find to the paragraph's own
# text in this host, so re-wrap the paragraph's Start/End into an explicit
# Document.Range before searching -- that keeps each Find call from
# wandering into other paragraphs that contain the same punctuation.
function Get-ParaRange($p) {
    $pr = $p.Range
    return $d.Range($pr.Start, $pr.End)
}

# --- Paragraph 1: {TEAM_LEADER} ({TEAM_LEADER_SERIAL}) -----------------
# -> {TEAM_LEADER} {TEAM_LEADER_SERIAL}
$rng = Get-ParaRange $paraTeamLeader
$rng.Find.Execute(" (", $true, $false, $false, $false, $false, $true, 0, $false, " ", 2)

$rng = Get-ParaRange $paraTeamLeader
$rng.Find.Execute(")", $true, $false, $false, $false, $false, $true, 0, $false, "", 2)

# --- Paragraph 2: {INSPECTOR} ({INSPECTOR_1_SERIAL})    ----------------
# -> {INSPECTOR} {INSPECTOR_1_SERIAL}   (trailing spaces removed too)
$rng = Get-ParaRange $paraInspector1
$rng.Find.Execute("(", $true, $false, $false, $false, $false, $true, 0, $false, "", 2)

$rng = Get-ParaRange $paraInspector1
$rng.Find.Execute(")   ", $true, $false, $false, $false, $false, $true, 0, $false, "", 2)

# --- Paragraph 3: {INSPECTOR_2} ({INSPECTOR_2_SERIAL})   ---------------
# -> {INSPECTOR_2} {INSPECTOR_2_SERIAL}   (trailing spaces kept)
$rng = Get-ParaRange $paraInspector2
$rng.Find.Execute("(", $true, $false, $false, $false, $false, $true, 0, $false, "", 2)

$rng = Get-ParaRange $paraInspector2
$rng.Find.Execute(")", $true, $false, $false, $false, $false, $true, 0, $false, "", 2)

# Re-seat the editing cursor mark ("_GoBack") right after the first
# character of INSPECTOR_2, the way Word leaves it behind after in-place
# retyping.
$rng = Get-ParaRange $paraInspector2
$found = $rng.Find.Execute("{INSPECTOR_2} {I", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
if ($found) {
    $bmRange = $d.Range($rng.End, $rng.End)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# --- Paragraph 4: {INSPECTOR_3} ({INSPECTOR_3_SERIAL})        ---------
# -> {INSPECTOR_3} {INSPECTOR_3_SERIAL}      (8 trailing spaces -> 6)
$rng = Get-ParaRange $paraInspector3
$rng.Find.Execute("(", $true, $false, $false, $false, $false, $true, 0, $false, "", 2)

$rng = Get-ParaRange $paraInspector3
$rng.Find.Execute(")", $true, $false, $false, $false, $false, $true, 0, $false, "", 2)

$rng = Get-ParaRange $paraInspector3
$rng.Find.Execute("        ", $true, $false, $false, $false, $false, $true, 0, $false, "      ", 2)
